# Update "想去人数" (F column) values on several sheets to match the
# regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 8503
$ws1.Range("F5").Value  = 8503
$ws1.Range("F6").Value  = 567
$ws1.Range("F7").Value  = 7533
$ws1.Range("F9").Value  = 632
$ws1.Range("F10").Value = 525
$ws1.Range("F15").Value = 12553
$ws1.Range("F18").Value = 2696
$ws1.Range("F19").Value = 5887
$ws1.Range("F22").Value = 3097
$ws1.Range("F24").Value = 142
$ws1.Range("F26").Value = 23
$ws1.Range("F31").Value = 1790
$ws1.Range("F33").Value = 159
$ws1.Range("F34").Value = 6231
$ws1.Range("F41").Value = 960

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 5

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 369
$ws3.Range("F3").Value = 528

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 369
$ws4.Range("F7").Value  = 528
$ws4.Range("F9").Value  = 8503
$ws4.Range("F10").Value = 567
$ws4.Range("F11").Value = 7533
$ws4.Range("F12").Value = 7533
$ws4.Range("F13").Value = 632
$ws4.Range("F14").Value = 525
$ws4.Range("F19").Value = 12553
$ws4.Range("F22").Value = 2696
$ws4.Range("F23").Value = 2696
$ws4.Range("F24").Value = 5887
$ws4.Range("F26").Value = 142
$ws4.Range("F28").Value = 23
$ws4.Range("F33").Value = 1790
$ws4.Range("F35").Value = 159
$ws4.Range("F36").Value = 6231
$ws4.Range("F45").Value = 960
